$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Keep only the "Sheet" tab; drop the three "Sec*" sheets so the
# remaining data requirements sheet is the sole/first tab.
[void]$wb.Worksheets.Item("Sec3").Delete()
[void]$wb.Worksheets.Item("Sec2").Delete()
[void]$wb.Worksheets.Item("Sec1").Delete()
